$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Riccardo Zaffoni "
$ws.Range("B11").Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Range("C11").Value = "Riccardo Zaffoni | U.S. Guarna"
$ws.Range("D11").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E11").Value = "Leonardo Viola | Shark Attack"
$ws.Range("F11").Value = "Daniele Ruzzenenti | Demobusters"
